$wb = $excel.ActiveWorkbook

# --- Stack sheet: add row 5 (entered first, 2021-07-21) ---
$wsStack = $wb.Worksheets.Item("Stack")
$wsStack.Activate()

$wsStack.Range("E5").Value = "DeleteMiddleElementOfAStack"
$wsStack.Range("B5").Value = "Delete middle element of a stack"
$wsStack.Range("C5").Value = "Medium"

$wsStack.Range("A4").Copy()
$wsStack.Range("A5").PasteSpecial(-4122)
$wsStack.Range("A5").Value = 44398

$wsStack.Application.CutCopyMode = $false
$wsStack.Range("A6").Select()

# --- General problems sheet: add row 8 (entered next, 2021-07-22) ---
$wsGeneral = $wb.Worksheets.Item("General problems")
$wsGeneral.Activate()

$wsGeneral.Range("E8").Value = "KthSymbolInGrammar"
$wsGeneral.Range("B8").Value = "Print Kth symbol in grammar"
$wsGeneral.Range("C8").Value = "Medium"

$wsGeneral.Range("A7").Copy()
$wsGeneral.Range("A8").PasteSpecial(-4122)
$wsGeneral.Range("A8").Value = 44399

$wsGeneral.Application.CutCopyMode = $false
$wsGeneral.Range("B8").Select()
